$wb = $excel.ActiveWorkbook

# --- "Test Results" sheet: mark Create + Read tests as passed (TRUE) for all data rows ---
$wsResults = $wb.Worksheets.Item("Test Results")
$wsResults.Range("B2:C24").Value = $true

# --- Update the remembered selection on each sheet (matches the saved UI state) ---
$wsDevice   = $wb.Worksheets.Item("Device")
$wsZone     = $wb.Worksheets.Item("Zone")
$wsCategory = $wb.Worksheets.Item("Category")

[void]$wsZone.Range("C32").Select()
[void]$wsCategory.Range("B12").Select()

# Device keeps its existing selection (E18); just touch it so ordering/tab state
# is driven entirely by the final Activate below.
[void]$wsDevice.Range("E18").Select()

# "Test Results" becomes the active sheet/tab, selection moves to G21.
[void]$wsResults.Range("G21").Select()
